# Test Menu part 2 - Update Data / Update Log -> Update Data, new Items, new Log rows.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Log" sheet: new Update-Data log rows (B5, B6, B8, B12).
# Written before the "Items" sheet edits so the new shared strings land in
# the same order as the authored workbook ("1s" before "Golden_Wind").
# ---------------------------------------------------------------------------
$wsLog = $wb.Worksheets.Item("Log")
$wsLog.Range("B5").Value = 1
$wsLog.Range("B6").Value = "1s"
$wsLog.Range("B8").Value = 11
$wsLog.Range("B12").Value = 1
[void]$wsLog.Range("B6:B7").Select()

# ---------------------------------------------------------------------------
# "Items" sheet: existing item (row 2) now has inventory, plus a brand new
# item (row 3) added via "Create New Item".
# ---------------------------------------------------------------------------
$wsItems = $wb.Worksheets.Item("Items")
$wsItems.Range("E2").Value = 1

$wsItems.Range("A3").Value = 13000006057
$wsItems.Range("B3").Value = "Golden_Wind"
$wsItems.Range("C3").Value = 50
$wsItems.Range("D3").Value = 15.25
$wsItems.Range("E3").Value = 1

# Column widths (best-fit) picked up for the newly populated columns.
$wsItems.Columns.Item(1).ColumnWidth = 11.166666666666666
$wsItems.Columns.Item(2).ColumnWidth = 12.451822916666666
$wsItems.Columns.Item(4).ColumnWidth = 5.166666666666667

# "Items" ends up the active sheet/cell after the Test Menu run.
[void]$wsItems.Range("E3").Select()
$wsItems.Activate()
